$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.2905003108896267
$ws.Cells.Item(2, 3).Value = 0.05300363562760424
$ws.Cells.Item(2, 4).Value = 0.07884844032159322
$ws.Cells.Item(2, 5).Value = 0.1643386443564268
$ws.Cells.Item(2, 7).Value = 0.5155318926930903
$ws.Cells.Item(2, 8).Value = 0.6533059411956543
$ws.Cells.Item(2, 11).Value = 0.267125269125728
$ws.Cells.Item(2, 13).Value = 0.2202148238570132
$ws.Cells.Item(2, 14).Value = 1.420034607490614
$ws.Cells.Item(2, 15).Value = 2.304186272167101
$ws.Cells.Item(3, 2).Value = 0.2573022933752043
$ws.Cells.Item(3, 3).Value = 0.04959231472099646
$ws.Cells.Item(3, 4).Value = 0.07148506819844158
$ws.Cells.Item(3, 5).Value = 0.1532663542963633
$ws.Cells.Item(3, 7).Value = 0.5144014372140475
$ws.Cells.Item(3, 8).Value = 0.656404492042256
$ws.Cells.Item(3, 11).Value = 0.2332290512571404
$ws.Cells.Item(3, 13).Value = 0.1979809654605091
$ws.Cells.Item(3, 14).Value = 1.435918653978847
$ws.Cells.Item(3, 15).Value = 2.308011874341858
$ws.Cells.Item(4, 2).Value = 0.2369331370846055
$ws.Cells.Item(4, 3).Value = 0.04748029964767397
$ws.Cells.Item(4, 4).Value = 0.06699663149181845
$ws.Cells.Item(4, 5).Value = 0.1465790906907145
$ws.Cells.Item(4, 7).Value = 0.5140529434752921
$ws.Cells.Item(4, 8).Value = 0.6585907717221886
$ws.Cells.Item(4, 11).Value = 0.2123733508069705
$ws.Cells.Item(4, 13).Value = 0.1844010535039899
$ws.Cells.Item(4, 14).Value = 1.446168424742691
$ws.Cells.Item(4, 15).Value = 2.311678139933875
$ws.Cells.Item(5, 2).Value = 0.2286366112731457
$ws.Cells.Item(5, 3).Value = 0.04661528658142799
$ws.Cells.Item(5, 4).Value = 0.06517581070593792
$ws.Cells.Item(5, 5).Value = 0.1438817384347502
$ws.Cells.Item(5, 7).Value = 0.5139977602228711
$ws.Cells.Item(5, 8).Value = 0.6595530973469721
$ws.Cells.Item(5, 11).Value = 0.203863992538345
$ws.Cells.Item(5, 13).Value = 0.1788851572464907
$ws.Cells.Item(5, 14).Value = 1.450470197681301
$ws.Cells.Item(5, 15).Value = 2.313503355255293
$ws.Cells.Item(6, 2).Value = 0.2272592374056046
$ws.Cells.Item(6, 3).Value = 0.04647139023148128
$ws.Cells.Item(6, 4).Value = 0.06487396418221181
$ws.Cells.Item(6, 5).Value = 0.1434355162949359
$ws.Cells.Item(6, 7).Value = 0.5139938392382533
$ws.Cells.Item(6, 8).Value = 0.6597172044617636
$ws.Cells.Item(6, 11).Value = 0.202450398045869
$ws.Cells.Item(6, 13).Value = 0.1779703352231721
$ws.Cells.Item(6, 14).Value = 1.451192046754282
$ws.Cells.Item(6, 15).Value = 2.313826431428126
$ws.Cells.Item(7, 2).Value = 0.2368212302120298
$ws.Cells.Item(7, 3).Value = 0.04746865134077893
$ws.Cells.Item(7, 4).Value = 0.06697204180719041
$ws.Cells.Item(7, 5).Value = 0.1465426011012738
$ws.Cells.Item(7, 7).Value = 0.5140518477813814
$ws.Cells.Item(7, 8).Value = 0.6586034608307614
$ws.Cells.Item(7, 11).Value = 0.2122586325817082
$ws.Cells.Item(7, 13).Value = 0.1843265912263021
$ws.Cells.Item(7, 14).Value = 1.446225934264705
$ws.Cells.Item(7, 15).Value = 2.311701414578863
$ws.Cells.Item(8, 2).Value = 0.2790508478228162
$ws.Cells.Item(8, 3).Value = 0.05183105593988557
$ws.Cells.Item(8, 4).Value = 0.07630277781188965
$ws.Cells.Item(8, 5).Value = 0.16049772695348
$ws.Cells.Item(8, 7).Value = 0.5150703326704473
$ws.Cells.Item(8, 8).Value = 0.6543154546481418
$ws.Cells.Item(8, 11).Value = 0.2554470631251604
$ws.Cells.Item(8, 13).Value = 0.2125336763846448
$ws.Cells.Item(8, 14).Value = 1.425408270991313
$ws.Cells.Item(8, 15).Value = 2.305231879449622
$ws.Cells.Item(9, 2).Value = 0.3619645948435561
$ws.Cells.Item(9, 3).Value = 0.06024595412694111
$ws.Cells.Item(9, 4).Value = 0.09485944867509488
$ws.Cells.Item(9, 5).Value = 0.1887559033010717
$ws.Cells.Item(9, 7).Value = 0.5198141786773363
$ws.Cells.Item(9, 8).Value = 0.6481564743262425
$ws.Cells.Item(9, 11).Value = 0.3397827085833285
$ws.Cells.Item(9, 13).Value = 0.2684208594187396
$ws.Cells.Item(9, 14).Value = 1.388528823392649
$ws.Cells.Item(9, 15).Value = 2.303003503440408
$ws.Cells.Item(10, 2).Value = 0.4229304981409143
$ws.Cells.Item(10, 3).Value = 0.06634197695834132
$ws.Cells.Item(10, 4).Value = 0.1086519344883499
$ws.Cells.Item(10, 5).Value = 0.2100773150175641
$ws.Cells.Item(10, 7).Value = 0.5249812073324733
$ws.Cells.Item(10, 8).Value = 0.6450011954734265
$ws.Cells.Item(10, 11).Value = 0.4015146870138153
$ws.Cells.Item(10, 13).Value = 0.309839962488482
$ws.Cells.Item(10, 14).Value = 1.363838576692491
$ws.Cells.Item(10, 15).Value = 2.30775417382344
$ws.Cells.Item(11, 2).Value = 0.450673787172974
$ws.Cells.Item(11, 3).Value = 0.06909623327130987
$ws.Cells.Item(11, 4).Value = 0.1149612697815456
$ws.Cells.Item(11, 5).Value = 0.219902211449245
$ws.Cells.Item(11, 7).Value = 0.527698648200527
$ws.Cells.Item(11, 8).Value = 0.6438628615397306
$ws.Cells.Item(11, 11).Value = 0.4295461604806405
$ws.Cells.Item(11, 13).Value = 0.3287627404166003
$ws.Cells.Item(11, 14).Value = 1.353128753974067
$ws.Cells.Item(11, 15).Value = 2.311305327957939
$ws.Cells.Item(12, 2).Value = 0.4611804903202597
$ws.Cells.Item(12, 3).Value = 0.07013645558623693
$ws.Cells.Item(12, 4).Value = 0.117355480715517
$ws.Cells.Item(12, 5).Value = 0.223640958706163
$ws.Cells.Item(12, 7).Value = 0.528780544605695
$ws.Cells.Item(12, 8).Value = 0.6434744818151472
$ws.Cells.Item(12, 11).Value = 0.4401533347431155
$ws.Cells.Item(12, 13).Value = 0.3359400385458073
$ws.Cells.Item(12, 14).Value = 1.349148322748579
$ws.Cells.Item(12, 15).Value = 2.312850110690619
$ws.Cells.Item(13, 2).Value = 0.4589176501106351
$ws.Cells.Item(13, 3).Value = 0.06991254830839466
$ws.Cells.Item(13, 4).Value = 0.1168396227016615
$ws.Cells.Item(13, 5).Value = 0.2228349367518803
$ws.Cells.Item(13, 7).Value = 0.5285451863504989
$ws.Cells.Item(13, 8).Value = 0.643556228462387
$ws.Cells.Item(13, 11).Value = 0.437869241127089
$ws.Cells.Item(13, 13).Value = 0.3343937606889185
$ws.Cells.Item(13, 14).Value = 1.350002236090035
$ws.Cells.Item(13, 15).Value = 2.312508514973189
$ws.Cells.Item(14, 2).Value = 0.4515381641160729
$ws.Cells.Item(14, 3).Value = 0.06918186844534091
$ws.Cells.Item(14, 4).Value = 0.1151581429253667
$ws.Cells.Item(14, 5).Value = 0.2202094330911208
$ws.Cells.Item(14, 7).Value = 0.5277865965182968
$ws.Cells.Item(14, 8).Value = 0.6438300540430504
$ws.Cells.Item(14, 11).Value = 0.430418976683228
$ws.Cells.Item(14, 13).Value = 0.3293529875044641
$ws.Cells.Item(14, 14).Value = 1.352799774886896
$ws.Cells.Item(14, 15).Value = 2.311428408251487
$ws.Cells.Item(15, 2).Value = 0.4470181223013583
$ws.Cells.Item(15, 3).Value = 0.06873394588878057
$ws.Cells.Item(15, 4).Value = 0.1141288383072379
$ws.Cells.Item(15, 5).Value = 0.2186036217912388
$ws.Cells.Item(15, 7).Value = 0.5273288251783157
$ws.Cells.Item(15, 8).Value = 0.6440033378247279
$ws.Cells.Item(15, 11).Value = 0.4258544539057141
$ws.Cells.Item(15, 13).Value = 0.3262668859678683
$ws.Cells.Item(15, 14).Value = 1.354523137395
$ws.Cells.Item(15, 15).Value = 2.310792866959304
$ws.Cells.Item(16, 2).Value = 0.4211175628793455
$ws.Cells.Item(16, 3).Value = 0.06616159670362265
$ws.Cells.Item(16, 4).Value = 0.1082403068050724
$ws.Cells.Item(16, 5).Value = 0.2094377809941435
$ws.Cells.Item(16, 7).Value = 0.5248110080374317
$ws.Cells.Item(16, 8).Value = 0.6450815626191257
$ws.Cells.Item(16, 11).Value = 0.3996817058823865
$ws.Cells.Item(16, 13).Value = 0.3086049449924388
$ws.Cells.Item(16, 14).Value = 1.364548999000234
$ws.Cells.Item(16, 15).Value = 2.30755008177033
$ws.Cells.Item(17, 2).Value = 0.4052305401632736
$ws.Cells.Item(17, 3).Value = 0.06457868363411023
$ws.Cells.Item(17, 4).Value = 0.1046368363463444
$ws.Cells.Item(17, 5).Value = 0.203847166195942
$ws.Cells.Item(17, 7).Value = 0.5233604577249622
$ws.Cells.Item(17, 8).Value = 0.6458190719142465
$ws.Cells.Item(17, 11).Value = 0.3836122681052245
$ws.Cells.Item(17, 13).Value = 0.2977906693640122
$ws.Cells.Item(17, 14).Value = 1.370833279267189
$ws.Cells.Item(17, 15).Value = 2.3059168704211
$ws.Cells.Item(18, 2).Value = 0.3960936753105955
$ws.Cells.Item(18, 3).Value = 0.06366646177032464
$ws.Cells.Item(18, 4).Value = 0.1025675182891206
$ws.Cells.Item(18, 5).Value = 0.2006434100473342
$ws.Cells.Item(18, 7).Value = 0.5225606713826494
$ws.Cells.Item(18, 8).Value = 0.646271226181014
$ws.Cells.Item(18, 11).Value = 0.3743648057063353
$ws.Cells.Item(18, 13).Value = 0.2915782099884083
$ws.Cells.Item(18, 14).Value = 1.374496948177215
$ws.Cells.Item(18, 15).Value = 2.305108319619791
$ws.Cells.Item(19, 2).Value = 0.393000263608684
$ws.Cells.Item(19, 3).Value = 0.06335729607931739
$ws.Cells.Item(19, 4).Value = 0.1018674512105662
$ws.Cells.Item(19, 5).Value = 0.1995606968606864
$ws.Cells.Item(19, 7).Value = 0.5222958054130089
$ws.Cells.Item(19, 8).Value = 0.646429120600061
$ws.Cells.Item(19, 11).Value = 0.3712329720829359
$ws.Cells.Item(19, 13).Value = 0.2894760876668414
$ws.Cells.Item(19, 14).Value = 1.375745836476372
$ws.Cells.Item(19, 15).Value = 2.304857023665846
$ws.Cells.Item(20, 2).Value = 0.4069216478735029
$ws.Cells.Item(20, 3).Value = 0.06474737115637197
$ws.Cells.Item(20, 4).Value = 0.1050200903812453
$ws.Cells.Item(20, 5).Value = 0.20444107216003
$ws.Cells.Item(20, 7).Value = 0.523511296669497
$ws.Cells.Item(20, 8).Value = 0.6457376694034593
$ws.Cells.Item(20, 11).Value = 0.3853233816810189
$ws.Cells.Item(20, 13).Value = 0.298941078192513
$ws.Cells.Item(20, 14).Value = 1.370159222759987
$ws.Cells.Item(20, 15).Value = 2.306077187135344
$ws.Cells.Item(21, 2).Value = 0.4537056766710634
$ws.Cells.Item(21, 3).Value = 0.06939656194229826
$ws.Cells.Item(21, 4).Value = 0.1156518988467496
$ws.Cells.Item(21, 5).Value = 0.2209801096670745
$ws.Cells.Item(21, 7).Value = 0.528007977359735
$ws.Cells.Item(21, 8).Value = 0.6437484666615205
$ws.Cells.Item(21, 11).Value = 0.4326075134495397
$ws.Cells.Item(21, 13).Value = 0.3308332685124142
$ws.Cells.Item(21, 14).Value = 1.35197602968562
$ws.Cells.Item(21, 15).Value = 2.311740231617847
$ws.Cells.Item(22, 2).Value = 0.4842869035138335
$ws.Cells.Item(22, 3).Value = 0.07241899983205258
$ws.Cells.Item(22, 4).Value = 0.1226295402230306
$ws.Cells.Item(22, 5).Value = 0.23189592364173
$ws.Cells.Item(22, 7).Value = 0.531254973787199
$ws.Cells.Item(22, 8).Value = 0.6426971818460743
$ws.Cells.Item(22, 11).Value = 0.4634651272395161
$ws.Cells.Item(22, 13).Value = 0.3517445803211032
$ws.Cells.Item(22, 14).Value = 1.340530318902338
$ws.Cells.Item(22, 15).Value = 2.31660740073454
$ws.Cells.Item(23, 2).Value = 0.4679648065274193
$ws.Cells.Item(23, 3).Value = 0.07080735323467025
$ws.Cells.Item(23, 4).Value = 0.1189027887005238
$ws.Cells.Item(23, 5).Value = 0.2260601287608992
$ws.Cells.Item(23, 7).Value = 0.5294937628967062
$ws.Cells.Item(23, 8).Value = 0.6432355187286021
$ws.Cells.Item(23, 11).Value = 0.447000129687126
$ws.Cells.Item(23, 13).Value = 0.3405776043284092
$ws.Cells.Item(23, 14).Value = 1.346598994745668
$ws.Cells.Item(23, 15).Value = 2.313902958608281
$ws.Cells.Item(24, 2).Value = 0.4061571080229669
$ws.Cells.Item(24, 3).Value = 0.0646711143278651
$ws.Cells.Item(24, 4).Value = 0.1048468138634462
$ws.Cells.Item(24, 5).Value = 0.2041725350122618
$ws.Cells.Item(24, 7).Value = 0.5234429960006395
$ws.Cells.Item(24, 8).Value = 0.6457743837885488
$ws.Cells.Item(24, 11).Value = 0.3845498150128606
$ws.Cells.Item(24, 13).Value = 0.2984209633763015
$ws.Cells.Item(24, 14).Value = 1.370463805619622
$ws.Cells.Item(24, 15).Value = 2.306004301755763
$ws.Cells.Item(25, 2).Value = 0.3395246709956723
$ws.Cells.Item(25, 3).Value = 0.0579845881327401
$ws.Cells.Item(25, 4).Value = 0.08981153735874159
$ws.Cells.Item(25, 5).Value = 0.1810141234389562
$ws.Cells.Item(25, 7).Value = 0.5182360414621741
$ws.Cells.Item(25, 8).Value = 0.6495819748121079
$ws.Cells.Item(25, 11).Value = 0.3170071071906762
$ws.Cells.Item(25, 13).Value = 0.2532395618572707
$ws.Cells.Item(25, 14).Value = 1.398083398149121
$ws.Cells.Item(25, 15).Value = 2.302485441086418
